$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells are not auto-converted to numbers/dates by Excel
# by forcing Text number format before assignment.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.426.23'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.873.58'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.018'
$ws.Range('E4').Value = '  +0.78%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.10'
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('E6').Value = '  +0.89%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5120'
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3961'
$ws.Range('E8').Value = '  +1.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08457'
$ws.Range('E9').Value = '  +0.76%  '
$ws.Range('E10').Value = '  -0.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.90'
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.263'
$ws.Range('E12').Value = '  +0.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.870.00'
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('E14').Value = '  -0.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.228'
$ws.Range('E15').Value = '  -0.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.018'
$ws.Range('E16').Value = '  +0.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001111'
$ws.Range('E17').Value = '  +0.64%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.94'
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('E19').Value = '  +1.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.72'
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.946'
$ws.Range('E22').Value = '  -1.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.445.03'
$ws.Range('E23').Value = '  +0.66%  '
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.280'
$ws.Range('E25').Value = '  +1.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.079.39'
$ws.Range('E26').Value = '  -0.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '161.89'
$ws.Range('E27').Value = '  +1.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.67'
$ws.Range('E28').Value = '  -0.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.374'
$ws.Range('E29').Value = '  -4.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.54'
$ws.Range('E30').Value = '  +0.74%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1054'
$ws.Range('E31').Value = '  -0.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.041'
$ws.Range('E32').Value = '  -0.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.789'
$ws.Range('E33').Value = '  -1.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.647'
$ws.Range('E34').Value = '  +0.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02436'
$ws.Range('E35').Value = '  -0.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06460'
$ws.Range('E36').Value = '  -1.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2183'
$ws.Range('E37').Value = '  -1.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.879'
$ws.Range('E38').Value = '  -7.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.260'
$ws.Range('E39').Value = '  +1.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.185'
$ws.Range('E40').Value = '  -1.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6386'
$ws.Range('E41').Value = '  -1.53%  '
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.24'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6018'
$ws.Range('E44').Value = '  -1.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.97'
$ws.Range('E45').Value = '  -1.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.712'
$ws.Range('E46').Value = '  +0.37%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.992'
$ws.Range('E47').Value = '  -1.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.205'
$ws.Range('E48').Value = '  -6.01%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '121.84'
$ws.Range('E49').Value = '  +0.64%  '
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.204'
$ws.Range('E50').Value = '  -2.67%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06850'
$ws.Range('E51').Value = '  -0.99%  '
